# edit.ps1 - apply the two changes described by the commit diff:
#   1. The table on slide 6 switches from table style {689212BB-...} to
#      {8DE539F3-...}.
#   2. The presentation's theme colour scheme (ppt/theme/theme1.xml - the
#      theme used by the slide master, i.e. every slide) changes from the
#      "Integral" palette to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($j = 1; $j -le $slide6.Shapes.Count; $j++) {
    $shape = $slide6.Shapes.Item($j)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{8DE539F3-40D5-4FCC-AF83-3C6A4E7C4FCF}")
    }
}

# --- 2. Theme colour scheme: Integral -> Office Theme ---------------------
# Helper that packs R,G,B (0-255 each) into the BGR-ordered long that the
# ThemeColor.RGB property expects (standard VBA/COM RGB() packing).
function ToRGB([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme values, in ThemeColorScheme.Item() order
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
$officeThemeColors = @(
    (ToRGB 0x00 0x00 0x00),  # dk1
    (ToRGB 0xFF 0xFF 0xFF),  # lt1
    (ToRGB 0x44 0x54 0x6A),  # dk2
    (ToRGB 0xE7 0xE6 0xE6),  # lt2
    (ToRGB 0x5B 0x9B 0xD5),  # accent1
    (ToRGB 0xED 0x7D 0x31),  # accent2
    (ToRGB 0xA5 0xA5 0xA5),  # accent3
    (ToRGB 0xFF 0xC0 0x00),  # accent4
    (ToRGB 0x44 0x72 0xC4),  # accent5
    (ToRGB 0x70 0xAD 0x47),  # accent6
    (ToRGB 0x05 0x63 0xC1),  # hlink
    (ToRGB 0x95 0x4F 0x72)   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
